# Apply the 09.11.2022 update to "CUENTAS DE CAJA IVSA.xlsx":
# the blank row-group that used to sit at A17:E19 is moved up to A4:E6 and
# filled in with the new recaudadora account data, while the old filled-in
# rows (15:16) and the vacated block (17:19) are cleared out.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring the formatting of the (previously blank) A17:E19 block up to A4:E6
# so the new rows pick up the right borders/number formats, mirroring what
# a cut-and-paste of that block would have done.
$ws.Range("A17:E19").Copy($ws.Range("A4"))

# --- New recaudadora accounts written into rows 4-6 ---
$ws.Range("B4").Value = "CAJA GENERAL M/N - RECAUDADORA AG HONDURAS CBBA"
$ws.Range("C4").Value = 101010112
$ws.Range("D4").Value = 101044022
$ws.Range("E4").Value = "BISA 0031"

$ws.Range("B5").Value = "CAJA GENERAL M/N - RECAUDADORA AG CALAMA CBBA"
$ws.Range("C5").Value = 101010113
$ws.Range("D5").Value = 101044022
$ws.Range("E5").Value = "BISA 0031"
$ws.Rows.Item(5).RowHeight = 15.75

$ws.Range("A6").Value = "ENVIO DE RECAUDACION L-M-V"
$ws.Range("B6").Value = "CAJA GENERAL M/N - RECAUDADORA CENTRAL SUCRE"
$ws.Range("C6").Value = 101010114
$ws.Range("D6").Value = 101020101
$ws.Range("E6").Value = "ETV M.N"

# --- Clear out the now-obsolete data previously shown in rows 15-19 ---
$ws.Range("B15:E15").ClearContents()
$ws.Range("A16:E16").ClearContents()
$ws.Range("A17:E17").Clear()
$ws.Range("A18:E18").Clear()
$ws.Range("A19:E19").Clear()

# --- Update the window's selection/scroll to match the saved view ---
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A4:E6").Select()
